$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.394.42"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").Value = "2.638.70"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.01"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.86"
$ws.Range("E6").Value = "  -1.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "2.667.44"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  +2.76%  "

$ws.Range("E11").Value = "  -0.90%  "

$ws.Range("E12").Value = "  -1.33%  "

$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "3.100.49"
$ws.Range("E14").Value = "  +0.43%  "

$ws.Range("D15").Value = "59.186.68"
$ws.Range("E15").Value = "  -2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").Value = "2.660.46"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.94"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.22"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.03"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  +1.06%  "

$ws.Range("D26").Value = "2.760.34"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.162"
$ws.Range("E28").Value = "  -2.30%  "

$ws.Range("D29").Value = "0.0₃0842"
$ws.Range("E29").Value = "  +0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.76"
$ws.Range("E31").Value = "  +11.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.997"
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.07"
$ws.Range("E33").Value = "  -0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.07"
$ws.Range("E35").Value = "  +19.92%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "149.47"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.08"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  -1.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.52"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  -0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "286.61"
$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.79"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0234"
$ws.Range("E49").Value = "  -1.16%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.80"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.25"
$ws.Range("E51").Value = "  -0.98%  "
